$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Logout" row (row 6) entirely, shifting rows below it up.
$ws.Rows(6).Delete()

# Mark the two Login Page test cases as Done in the Status column.
$ws.Range("D4").Value = "Done"
$ws.Range("D5").Value = "Done"

# Leave the selection where the user last clicked.
[void]$ws.Range("D6").Select()
